$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B34").Value = 253
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 403
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 1
